$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily row: day 45749 (2025-04-02) with Chase/Bryce/Zach win counts
$ws.Range("A9").Value = 45749
$ws.Range("B9").Value = 34
$ws.Range("C9").Value = 27
$ws.Range("D9").Value = 34

# The "last row" date-only format moves from A8 to the new last row A9;
# A8 reverts to the regular row's date/time format.
$ws.Range("A8").NumberFormat = $ws.Range("A7").NumberFormat
$ws.Range("A9").NumberFormat = "YYYY-MM-DD"
